$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.016355905698264
$ws.Cells.Item(2, 4).Value = 1.022479483613737
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.014699180095862
$ws.Cells.Item(2, 9).Value = 1.026597905994943
$ws.Cells.Item(2, 10).Value = 1.02157618415187
$ws.Cells.Item(2, 11).Value = 1.025313742829508
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.017556549926796
$ws.Cells.Item(2, 14).Value = 1.023026939509792

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.017389262125247
$ws.Cells.Item(3, 4).Value = 1.023228818503399
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.016376763410981
$ws.Cells.Item(3, 9).Value = 1.026739901939617
$ws.Cells.Item(3, 10).Value = 1.02224451437671
$ws.Cells.Item(3, 11).Value = 1.025869909336335
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.019036714724625
$ws.Cells.Item(3, 14).Value = 1.02369621884021

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.018057337339677
$ws.Cells.Item(4, 4).Value = 1.023712977436805
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.017461496354829
$ws.Cells.Item(4, 9).Value = 1.02683012277388
$ws.Cells.Item(4, 10).Value = 1.022675861717926
$ws.Cells.Item(4, 11).Value = 1.026228417103951
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.019993265129218
$ws.Cells.Item(4, 14).Value = 1.024128178744127

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.018338060112591
$ws.Cells.Item(5, 4).Value = 1.023916348077487
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.017917340599199
$ws.Cells.Item(5, 9).Value = 1.026867654294941
$ws.Cells.Item(5, 10).Value = 1.022856936295487
$ws.Cells.Item(5, 11).Value = 1.026378806620884
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.020395115144845
$ws.Cells.Item(5, 14).Value = 1.024309510468357

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.018385186773475
$ws.Cells.Item(6, 4).Value = 1.023950484958692
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.017993868683567
$ws.Cells.Item(6, 9).Value = 1.026873932706341
$ws.Cells.Item(6, 10).Value = 1.022887324075736
$ws.Cells.Item(6, 11).Value = 1.026404038516425
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.020462571097077
$ws.Cells.Item(6, 14).Value = 1.02433994140274

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.018061088903551
$ws.Cells.Item(7, 4).Value = 1.023715695551894
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.017467588055177
$ws.Cells.Item(7, 9).Value = 1.026830625832979
$ws.Cells.Item(7, 10).Value = 1.02267828228059
$ws.Cells.Item(7, 11).Value = 1.026230427900756
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.019998635775003
$ws.Cells.Item(7, 14).Value = 1.024130602744267

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.01670525354399
$ws.Cells.Item(8, 4).Value = 1.022732872281074
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.01526629278056
$ws.Cells.Item(8, 9).Value = 1.026646237971952
$ws.Cells.Item(8, 10).Value = 1.021802279403886
$ws.Cells.Item(8, 11).Value = 1.025501985779147
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.018057034474654
$ws.Cells.Item(8, 14).Value = 1.023253355842999

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.014311612582236
$ws.Cells.Item(9, 4).Value = 1.020995548613844
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.011381020659848
$ws.Cells.Item(9, 9).Value = 1.026308604200629
$ws.Cells.Item(9, 10).Value = 1.020250118624637
$ws.Cells.Item(9, 11).Value = 1.024207862492533
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.014626057207372
$ws.Cells.Item(9, 14).Value = 1.021698990817409

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.012712723436749
$ws.Cells.Item(10, 4).Value = 1.019833629035883
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.008786081658748
$ws.Cells.Item(10, 9).Value = 1.026074955273946
$ws.Cells.Item(10, 10).Value = 1.01920953114614
$ws.Cells.Item(10, 11).Value = 1.023338000285984
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.012331815916494
$ws.Cells.Item(10, 14).Value = 1.020656925585335

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.01201961876227
$ws.Cells.Item(11, 4).Value = 1.019329618672648
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.00766120077835
$ws.Cells.Item(11, 9).Value = 1.025971750119433
$ws.Cells.Item(11, 10).Value = 1.018757548354732
$ws.Cells.Item(11, 11).Value = 1.022959642601455
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.011336641464614
$ws.Cells.Item(11, 14).Value = 1.020204300926519

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.011762050023323
$ws.Cells.Item(12, 4).Value = 1.019142272046536
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.007243171375214
$ws.Cells.Item(12, 9).Value = 1.025933109332702
$ws.Cells.Item(12, 10).Value = 1.018589449607373
$ws.Cells.Item(12, 11).Value = 1.022818847050627
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.010966717286365
$ws.Cells.Item(12, 14).Value = 1.020035963459659

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.011817304806159
$ws.Cells.Item(13, 4).Value = 1.019182464643873
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.007332849247436
$ws.Cells.Item(13, 9).Value = 1.025941411754956
$ws.Cells.Item(13, 10).Value = 1.018625516989092
$ws.Cells.Item(13, 11).Value = 1.022849059805086
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.011046079731223
$ws.Cells.Item(13, 14).Value = 1.020072082061196

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.011998330487152
$ws.Cells.Item(14, 4).Value = 1.01931413529534
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.007626650446014
$ws.Cells.Item(14, 9).Value = 1.025968562299689
$ws.Cells.Item(14, 10).Value = 1.018743657598536
$ws.Cells.Item(14, 11).Value = 1.022948009639017
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.011306069017652
$ws.Cells.Item(14, 14).Value = 1.020190390443857

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.012109850570729
$ws.Cells.Item(15, 4).Value = 1.019395244051408
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.007807644450331
$ws.Cells.Item(15, 9).Value = 1.025985250117173
$ws.Cells.Item(15, 10).Value = 1.018816419757114
$ws.Cells.Item(15, 11).Value = 1.023008941856873
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.011466220682681
$ws.Cells.Item(15, 14).Value = 1.020263255933046

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.01275870598034
$ws.Cells.Item(16, 4).Value = 1.019867059693041
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.00886070900245
$ws.Cells.Item(16, 9).Value = 1.02608176179127
$ws.Cells.Item(16, 10).Value = 1.019239498086258
$ws.Cells.Item(16, 11).Value = 1.023363074702396
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.01239782469306
$ws.Cells.Item(16, 14).Value = 1.020686935081945

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.013165506661508
$ws.Cells.Item(17, 4).Value = 1.020162777939031
$ws.Cells.Item(17, 5).Value = 0.9894763578477731
$ws.Cells.Item(17, 6).Value = 1.009520925502114
$ws.Cells.Item(17, 9).Value = 1.026141756367808
$ws.Cells.Item(17, 10).Value = 1.01950450755309
$ws.Cells.Item(17, 11).Value = 1.023584756601122
$ws.Cells.Item(17, 12).Value = 0.9930127773692701
$ws.Cells.Item(17, 13).Value = 1.012981719845379
$ws.Cells.Item(17, 14).Value = 1.020952320892617

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.013402711989008
$ws.Cells.Item(18, 4).Value = 1.020335179426046
$ws.Cells.Item(18, 5).Value = 0.9897087662937551
$ws.Cells.Item(18, 6).Value = 1.009905898189045
$ws.Cells.Item(18, 9).Value = 1.026176553992585
$ws.Cells.Item(18, 10).Value = 1.019658947975695
$ws.Cells.Item(18, 11).Value = 1.023713895791683
$ws.Cells.Item(18, 12).Value = 0.9932001317071766
$ws.Cells.Item(18, 13).Value = 1.013322127479442
$ws.Cells.Item(18, 14).Value = 1.021106980638337

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.013483580287965
$ws.Cells.Item(19, 4).Value = 1.020393949320443
$ws.Cells.Item(19, 5).Value = 0.9897880325774039
$ws.Cells.Item(19, 6).Value = 1.010037143800917
$ws.Cells.Item(19, 9).Value = 1.026188385812848
$ws.Cells.Item(19, 10).Value = 1.019711585299954
$ws.Cells.Item(19, 11).Value = 1.02375790111249
$ws.Cells.Item(19, 12).Value = 0.993264023964098
$ws.Cells.Item(19, 13).Value = 1.013438169363255
$ws.Cells.Item(19, 14).Value = 1.021159692713633

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.013121868535433
$ws.Cells.Item(20, 4).Value = 1.020131059067889
$ws.Cells.Item(20, 5).Value = 0.9894336180360677
$ws.Cells.Item(20, 6).Value = 1.009450103050661
$ws.Cells.Item(20, 9).Value = 1.026135339813285
$ws.Cells.Item(20, 10).Value = 1.019476088533656
$ws.Cells.Item(20, 11).Value = 1.023560989214044
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.012919090918389
$ws.Cells.Item(20, 14).Value = 1.020923861514916

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.011945026234539
$ws.Cells.Item(21, 4).Value = 1.019275365296149
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.007540138877542
$ws.Cells.Item(21, 9).Value = 1.025960575584653
$ws.Cells.Item(21, 10).Value = 1.018708874005103
$ws.Cells.Item(21, 11).Value = 1.022918878472484
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.011229516233937
$ws.Cells.Item(21, 14).Value = 1.020155557453732

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.01120441023922
$ws.Cells.Item(22, 4).Value = 1.018736576709593
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.006338114802934
$ws.Cells.Item(22, 9).Value = 1.025848924992503
$ws.Cells.Item(22, 10).Value = 1.018225267145674
$ws.Cells.Item(22, 11).Value = 1.02251367249219
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.010165635051638
$ws.Cells.Item(22, 14).Value = 1.019671263817091

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.011597090846267
$ws.Cells.Item(23, 4).Value = 1.019022272934744
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.006975443226797
$ws.Cells.Item(23, 9).Value = 1.025908280892351
$ws.Cells.Item(23, 10).Value = 1.01848175330769
$ws.Cells.Item(23, 11).Value = 1.022728620997495
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.010729771153859
$ws.Cells.Item(23, 14).Value = 1.019928114218876

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.013141586956791
$ws.Cells.Item(24, 4).Value = 1.020145391726139
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.009482105039209
$ws.Cells.Item(24, 9).Value = 1.026138239784325
$ws.Cells.Item(24, 10).Value = 1.019488930282719
$ws.Cells.Item(24, 11).Value = 1.023571729180314
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.012947390754662
$ws.Cells.Item(24, 14).Value = 1.020936721500735

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.014930969254094
$ws.Cells.Item(25, 4).Value = 1.021445338945091
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.012386257174939
$ws.Cells.Item(25, 9).Value = 1.026397398702084
$ws.Cells.Item(25, 10).Value = 1.020652408352513
$ws.Cells.Item(25, 11).Value = 1.024543674910681
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.015514234164774
$ws.Cells.Item(25, 14).Value = 1.022101851842842

